# Update Linux-projects, and small fixes problems
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# "Корректное отображение окон, масштабирование" (row 31) is now fully working.
$ws.Range("C31").Interior.Color = 5287936
$ws.Range("C31").Value = 1

# "Переключение языков" (row 32) is now fully working too; the old
# "partially working / language-code lookup broken" comment no longer applies.
$ws.Range("C32").Interior.Color = 5287936
$ws.Range("C32").Value = 1
$ws.Range("D32").ClearContents() | Out-Null

# Move the selection / viewport the way the author left it.
$ws.Range("D18").Select() | Out-Null
